$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-20 down to 12-21.
$ws.Rows(11).Insert()
$ws.Rows(11).RowHeight = 20

# Populate the new row 11 with the "測試環境" task entry.
$ws.Range("A11").Value = "測試環境"
$ws.Range("B11").Value = "[第一次]切轉演練"
$ws.Range("C11").Value = "第一次"
$ws.Range("D11").Value = 46174
$ws.Range("E11").Value = 46178
$ws.Range("F11").Value = "OO"
$ws.Range("G11").Value = "進行中"
$ws.Range("H11").Value = 46136
$ws.Range("J11").Value = 46034
$ws.Range("K11").Value = 46034
$ws.Range("L11").Value = "[第一次]切轉演練"
$ws.Range("M11").Value = "v"
$ws.Range("N11").Value = "v"
$ws.Range("O11").Value = "v"
$ws.Range("P11").Value = "html://127.0.0.1"

# Update sheet view: zoom to 115% and select D11 (matches reviewed cell position).
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("D11").Select()
